# Update Name of Algo
# Apply updated imputed values to the result_data_KNN worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B11").Value = 6.782000000000001
$ws.Range("A12").Value = -21.489

$ws.Range("B23").Value = 8.162000000000001

$ws.Range("B28").Value = 5.196000000000001

$ws.Range("A32").Value = -21.23
$ws.Range("B32").Value = 7.417

$ws.Range("B34").Value = 7.119999999999999

$ws.Range("A36").Value = -20.724

$ws.Range("A38").Value = -20.331

$ws.Range("B42").Value = 9.103999999999999

$ws.Range("A46").Value = -21.697

$ws.Range("A54").Value = -21.066
$ws.Range("B54").Value = 5.524

$ws.Range("A55").Value = -22.184

$ws.Range("A67").Value = -21.422

$ws.Range("A69").Value = -21.422

$ws.Range("A72").Value = -21.621

$ws.Range("A91").Value = -20.847

$ws.Range("B97").Value = 5.191

$ws.Range("A99").Value = -21.326
$ws.Range("B99").Value = 5.483000000000001

$ws.Range("B101").Value = 5.217000000000001

$ws.Range("A104").Value = -21.437
